$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("F2").Value = 1.65
$ws.Range("G2").Value = 1.67
$ws.Range("H2").Value = 5.2
$ws.Range("I2").Value = 5.4
$ws.Range("J2").Value = 4.8
$ws.Range("L2").Value = 1.23
$ws.Range("N2").Value = 8.199999999999999
$ws.Range("Q2").Value = 1.41
$ws.Range("R2").Value = 1.92
$ws.Range("S2").Value = 2.04
$ws.Range("U2").Value = 2.9
$ws.Range("V2").Value = 1.23
$ws.Range("W2").Value = 2.48
$ws.Range("AC2").Value = 12.5
$ws.Range("AD2").Value = 21
$ws.Range("AE2").Value = 50
$ws.Range("AJ2").Value = 18.5
$ws.Range("AL2").Value = 21
$ws.Range("AM2").Value = 50
$ws.Range("AN2").Value = 5.3
$ws.Range("AO2").Value = 34

# Row 3
$ws.Range("F3").Value = 1.97
$ws.Range("G3").Value = 1.99
$ws.Range("H3").Value = 3.95
$ws.Range("I3").Value = 4.1
$ws.Range("N3").Value = 5.6
$ws.Range("P3").Value = 2.5
$ws.Range("T3").Value = 1.6
$ws.Range("U3").Value = 2.58
$ws.Range("V3").Value = 1.32
$ws.Range("W3").Value = 2
$ws.Range("X3").Value = 24
$ws.Range("AA3").Value = 75
$ws.Range("AJ3").Value = 23
$ws.Range("AN3").Value = 9.199999999999999
$ws.Range("AO3").Value = 30

# Row 4
$ws.Range("F4").Value = 1.77
$ws.Range("G4").Value = 1.93
$ws.Range("H4").Value = 4
$ws.Range("I4").Value = 4.7
$ws.Range("J4").Value = 4
$ws.Range("K4").Value = 4.7
$ws.Range("L4").Value = 1.31
$ws.Range("N4").Value = 4.8
$ws.Range("O4").Value = 1.21
$ws.Range("P4").Value = 2.36
$ws.Range("Q4").Value = 1.6
$ws.Range("R4").Value = 1.54
$ws.Range("S4").Value = 2.52
$ws.Range("T4").Value = 1.62
$ws.Range("U4").Value = 2.28
$ws.Range("W4").Value = 2.06
$ws.Range("Y4").Value = 44
$ws.Range("Z4").Value = 90
$ws.Range("AB4").Value = 12.5
$ws.Range("AC4").Value = 11
$ws.Range("AH4").Value = 25
$ws.Range("AN4").Value = 29

# Row 5
$ws.Range("F5").Value = 2.3
$ws.Range("G5").Value = 2.58
$ws.Range("H5").Value = 3.3
$ws.Range("I5").Value = 3.8
$ws.Range("J5").Value = 3.1
$ws.Range("L5").Value = 1.48
$ws.Range("M5").Value = 1.1
$ws.Range("N5").Value = 2.86
$ws.Range("O5").Value = 1.39
$ws.Range("P5").Value = 1.69
$ws.Range("Q5").Value = 2.22
$ws.Range("R5").Value = 1.25
$ws.Range("S5").Value = 4.1
$ws.Range("T5").Value = 1.88
$ws.Range("U5").Value = 1.93
$ws.Range("V5").Value = 1.37
$ws.Range("W5").Value = 1.63
$ws.Range("X5").Value = 11.5
$ws.Range("Y5").Value = 12
$ws.Range("Z5").Value = 25
$ws.Range("AA5").Value = 270
$ws.Range("AB5").Value = 9.199999999999999
$ws.Range("AC5").Value = 7.8
$ws.Range("AD5").Value = 15.5
$ws.Range("AE5").Value = 120
$ws.Range("AF5").Value = 15.5
$ws.Range("AG5").Value = 12.5
$ws.Range("AH5").Value = 19.5
$ws.Range("AI5").Value = 160
$ws.Range("AJ5").Value = 36
$ws.Range("AK5").Value = 30
$ws.Range("AL5").Value = 50
$ws.Range("AN5").Value = 28
$ws.Range("AO5").Value = 55

# Row 6
$ws.Range("F6").Value = 1.53
$ws.Range("G6").Value = 1.61
$ws.Range("J6").Value = 4
$ws.Range("K6").Value = 5
$ws.Range("L6").Value = 1.39
$ws.Range("M6").Value = 1.06
$ws.Range("N6").Value = 3.6
$ws.Range("O6").Value = 1.31
$ws.Range("P6").Value = 1.89
$ws.Range("Q6").Value = 1.9
$ws.Range("R6").Value = 1.34
$ws.Range("S6").Value = 3.35
$ws.Range("T6").Value = 1.99
$ws.Range("U6").Value = 1.81
$ws.Range("W6").Value = 2.58
$ws.Range("AB6").Value = 29
$ws.Range("AC6").Value = 10
$ws.Range("AF6").Value = 40
$ws.Range("AG6").Value = 19.5
$ws.Range("AJ6").Value = 180

# Row 7
$ws.Range("F7").Value = 13.5
$ws.Range("G7").Value = 21
$ws.Range("I7").Value = 1.29
$ws.Range("J7").Value = 6.6

# Row 8
$ws.Range("G8").Value = 5.2
$ws.Range("H8").Value = 1.84
$ws.Range("I8").Value = 4.8
$ws.Range("J8").Value = 3
$ws.Range("N8").Value = 1.1
$ws.Range("P8").Value = 1.67
$ws.Range("R8").Value = 1.22
$ws.Range("S8").Value = 1.05
$ws.Range("W8").Value = 1.24
$ws.Range("AC8").Value = 100

# Row 9
$ws.Range("O9").Value = 1.37
$ws.Range("AA9").Value = 65

# Row 10
$ws.Range("F10").Value = 1.56
$ws.Range("G10").Value = 1.57
$ws.Range("P10").Value = 2.48
$ws.Range("Q10").Value = 1.65
$ws.Range("R10").Value = 1.56
$ws.Range("W10").Value = 2.76
$ws.Range("Y10").Value = 28
$ws.Range("AO10").Value = 85

# Row 11
$ws.Range("G11").Value = 2.04
$ws.Range("I11").Value = 4.2
$ws.Range("L11").Value = 1.2
$ws.Range("N11").Value = 5.1
$ws.Range("O11").Value = 1.19
$ws.Range("P11").Value = 2.44
$ws.Range("Q11").Value = 1.57
$ws.Range("R11").Value = 1.59
$ws.Range("S11").Value = 2.4
$ws.Range("T11").Value = 1.56
$ws.Range("U11").Value = 2.4
$ws.Range("V11").Value = 1.31
$ws.Range("W11").Value = 1.96
$ws.Range("X11").Value = 26
$ws.Range("Y11").Value = 21
$ws.Range("AA11").Value = 75
$ws.Range("AB11").Value = 16
$ws.Range("AC11").Value = 10.5
$ws.Range("AD11").Value = 19
$ws.Range("AE11").Value = 42
$ws.Range("AF11").Value = 15.5
$ws.Range("AH11").Value = 18.5
$ws.Range("AJ11").Value = 24
$ws.Range("AK11").Value = 21
$ws.Range("AL11").Value = 29
$ws.Range("AM11").Value = 200
$ws.Range("AN11").Value = 10

# Row 12
$ws.Range("F12").Value = 2.34
$ws.Range("G12").Value = 2.36
$ws.Range("L12").Value = 1.41
$ws.Range("N12").Value = 3.95
$ws.Range("O12").Value = 1.32
$ws.Range("R12").Value = 1.4
$ws.Range("W12").Value = 1.73
$ws.Range("AF12").Value = 14.5
$ws.Range("AO12").Value = 32

# Row 13
$ws.Range("F13").Value = 8.800000000000001
$ws.Range("G13").Value = 9
$ws.Range("H13").Value = 1.42
$ws.Range("I13").Value = 1.43
$ws.Range("J13").Value = 5.2
$ws.Range("K13").Value = 5.3
$ws.Range("N13").Value = 5.1
$ws.Range("R13").Value = 1.55
$ws.Range("V13").Value = 3.3
$ws.Range("W13").Value = 1.12
$ws.Range("X13").Value = 22
$ws.Range("Y13").Value = 9.199999999999999
$ws.Range("Z13").Value = 8.6
$ws.Range("AE13").Value = 14
$ws.Range("AH13").Value = 24
$ws.Range("AJ13").Value = 270
$ws.Range("AN13").Value = 130
$ws.Range("AO13").Value = 5.8
